$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 23 de Marzo de 2020 a las 21:16"

$data = @(
  ,("China",81093,39,72703,5120,1749,9,3270)
  ,("Italia",63927,4789,7432,50418,3204,601,6077)
  ,("Estados Unidos",42163,8617,295,41356,1040,93,512)
  ,("España",33089,4321,3355,27527,2355,435,2207)
  ,("Alemania",29056,4183,422,28516,23,24,118)
  ,("Iran",23049,1411,8376,12861,0,127,1812)
  ,("Francia",19856,3838,2200,16796,2082,186,860)
  ,("Corea del Sur",8961,64,3166,5684,59,7,111)
  ,("Suiza",8547,1073,131,8298,141,20,118)
  ,("Reino Unido",6650,967,135,6180,20,54,335)
  ,("Paises Bajos",4749,545,2,4534,435,34,213)
  ,("Austria",4468,886,9,4438,14,5,21)
  ,("Belgica",3743,342,401,3254,322,13,88)
  ,("Noruega",2547,162,6,2531,41,3,10)
  ,("Portugal",2060,460,14,2023,47,9,23)
  ,("Suecia",2046,112,16,2005,104,4,25)
  ,("Canada",2035,565,320,1692,1,3,23)
  ,("Australia",1717,108,88,1622,11,0,7)
  ,("Brasil",1629,83,2,1602,18,0,25)
  ,("Malasia",1518,212,159,1345,57,4,14)
  ,("Dinamarca",1450,55,1,1425,55,11,24)
  ,("Israel",1442,371,41,1400,29,0,1)
  ,("Chequia",1236,116,7,1228,19,0,1)
  ,("Turquia",1236,0,0,1206,0,0,30)
  ,("Japon",1128,27,235,851,49,1,42)
  ,("Irlanda",1125,219,5,1114,29,2,6)
  ,("Ecuador",981,192,3,960,2,4,18)
  ,("Luxemburgo",875,77,6,861,3,0,8)
  ,("Pakistan",873,97,13,854,0,1,6)
  ,("Polonia",749,115,13,728,3,1,8)
  ,("Chile",746,114,11,733,7,1,2)
  ,("Tailandia",721,122,52,668,7,0,1)
  ,("Crucero",712,0,567,137,15,0,8)
  ,("Finlandia",700,74,10,689,13,0,1)
  ,("Grecia",695,71,29,649,35,2,17)
  ,("Islandia",588,20,51,536,14,0,1)
  ,("Indonesia",579,65,30,500,0,1,49)
  ,("Rumania",576,143,73,496,15,4,7)
  ,("Arabia Saudita",562,51,19,543,0,0,0)
  ,("Singapur",509,54,152,355,14,0,2)
  ,("Catar",501,7,33,468,6,0,0)
  ,("India",499,103,34,455,0,3,10)
  ,("Filipinas",462,82,18,411,1,8,33)
  ,("Eslovenia",442,28,0,439,12,1,3)
  ,("Rusia",438,71,17,420,0,0,1)
  ,("Sudafrica",402,128,4,398,0,0,0)
  ,("Peru",395,32,1,389,19,0,5)
  ,("Barein",377,43,164,211,3,0,2)
  ,("Egipto",366,39,68,279,0,5,19)
  ,("Hong Kong",356,38,100,252,4,0,4)
  ,("Estonia",352,26,4,348,4,0,0)
  ,("Mexico",316,65,4,310,1,0,2)
  ,("Croacia",315,61,5,309,5,0,1)
  ,("Panama",313,0,1,309,7,0,3)
  ,("Colombia",277,46,3,271,0,1,3)
  ,("Libano",267,19,8,255,4,0,4)
  ,("Argentina",266,0,27,235,0,0,4)
  ,("Irak",266,33,62,181,0,3,23)
  ,("Serbia",249,27,3,244,12,0,2)
  ,("Republica Dominicana",245,43,0,242,0,0,3)
  ,("Argelia",230,29,65,148,0,0,17)
  ,("Bulgaria",201,14,3,195,8,0,3)
  ,("Emiratos Arabes Unidos",198,45,41,155,2,0,2)
  ,("Taiwan",195,26,28,165,0,0,2)
  ,("Armenia",194,0,2,192,6,0,0)
  ,("Kuwait",189,1,30,159,5,0,0)
  ,("San Marino",187,12,4,163,13,0,20)
  ,("Eslovaquia",186,1,7,179,2,0,0)
  ,("Letonia",180,41,1,179,0,0,0)
  ,("Hungria",167,36,16,144,6,1,7)
  ,("Lituania",160,17,1,158,1,0,1)
  ,("Uruguay",158,0,0,158,2,0,0)
  ,("Costa Rica",158,24,2,154,2,0,2)
  ,("Marruecos",143,28,5,134,1,0,4)
  ,("Republica de Macedonia",136,21,1,133,1,1,2)
  ,("Principado de Andorra",133,20,1,131,2,0,1)
  ,("Bosnia y Herzegovina",132,6,2,129,1,0,1)
  ,("Jordania",127,15,1,126,0,0,0)
  ,("Vietnam",123,10,17,106,2,0,0)
  ,("Islas Feroe",118,3,23,95,0,0,0)
  ,("Republica de Chipre",116,21,3,112,3,0,1)
  ,("Moldavia",109,15,2,106,3,0,1)
  ,("Malta",107,17,2,105,1,0,0)
  ,("Albania",104,15,2,98,2,2,4)
  ,("Nueva Zelanda",102,0,0,102,0,0,0)
  ,("Burkina Faso",99,24,5,90,0,0,4)
  ,("Sri Lanka",97,15,2,95,2,0,0)
  ,("Brunei",91,3,2,89,2,0,0)
  ,("Tunez",89,14,1,85,11,0,3)
  ,("Camboya",87,3,2,85,1,0,0)
  ,("Bielorrusia",81,5,22,59,0,0,0)
  ,("Senegal",79,12,8,71,0,0,0)
  ,("Venezuela",77,7,15,62,2,0,0)
  ,("Ucrania",73,0,1,69,0,0,3)
  ,("Azerbaiyan",72,7,10,61,3,0,1)
  ,("Reunion",71,7,1,70,0,0,0)
  ,("Oman",66,11,17,49,0,0,0)
  ,("Kazajistan",62,2,0,62,0,0,0)
  ,("Guadalupe",62,4,0,61,4,0,1)
  ,("Georgia",61,7,8,53,1,0,0)
  ,("Estado de Palestina",59,0,17,42,0,0,0)
  ,("Camerun",56,16,2,54,0,0,0)
  ,("Martinica",53,9,0,52,7,0,1)
  ,("Liechtenstein",51,14,0,51,0,0,0)
  ,("Trinidad yTobago",51,1,0,51,0,0,0)
  ,("Uzbekistan",46,3,0,46,0,0,0)
  ,("Cuba",40,5,0,39,3,0,1)
  ,("Afganistan",40,0,1,38,0,0,1)
  ,("Consejo Danes para los Refugiados",36,6,0,35,0,0,1)
  ,("Mauricio",36,8,0,34,1,0,2)
  ,("Nigeria",36,6,2,33,0,1,1)
  ,("Banglades",33,6,5,25,0,1,3)
  ,("Puerto Rico",31,8,1,28,0,1,2)
  ,("Guam",29,2,0,28,0,0,1)
  ,("Honduras",27,1,0,27,0,0,0)
  ,("Bolivia",27,3,0,27,0,0,0)
  ,("Montenegro",27,6,0,26,0,1,1)
  ,("Costa de Marfil",25,11,2,23,0,0,0)
  ,("Macao",25,3,10,15,0,0,0)
  ,("Mayotte",24,13,0,24,0,0,0)
  ,("Ghana",24,1,0,23,0,0,1)
  ,("Monaco",23,0,1,22,0,0,0)
  ,("Paraguay",22,0,0,21,1,0,1)
  ,("Guyana",20,1,0,19,0,0,1)
  ,("Guatemala",20,1,0,19,0,0,1)
  ,("Guayana Francesa",20,2,6,14,0,0,0)
  ,("Ruanda",19,0,0,19,0,0,0)
  ,("Jamaica",19,0,2,16,0,0,1)
  ,("Polinesia Francesa",18,0,0,18,0,0,0)
  ,("Togo",18,2,0,18,0,0,0)
  ,("Barbados",17,3,0,17,0,0,0)
  ,("Islas Virgenes de los Estados Unidos",17,11,0,17,0,0,0)
  ,("Kirguistan",16,2,0,16,0,0,0)
  ,("Kenia",16,1,0,16,0,0,0)
  ,("Gibraltar",15,0,5,10,0,0,0)
  ,("Isla de Man",13,8,0,13,0,0,0)
  ,("Maldivas",13,0,5,8,0,0,0)
  ,("Tanzania",12,0,0,12,0,0,0)
  ,("Madagascar",12,9,0,12,0,0,0)
  ,("Etiopia",11,0,0,11,0,0,0)
  ,("Mongolia",10,0,0,10,0,0,0)
  ,("Guinea Ecuatorial",9,3,0,9,0,0,0)
  ,("Aruba",9,0,1,8,0,0,0)
  ,("Nueva Caledonia",8,4,0,8,0,0,0)
  ,("San Martin (Parte Francesa)",8,3,0,8,0,0,0)
  ,("Seychelles",7,0,0,7,0,0,0)
  ,("Haiti",6,4,0,6,0,0,0)
  ,("Bermudas",6,0,0,6,0,0,0)
  ,("Surinam",5,0,0,5,0,0,0)
  ,("Gabon",5,0,0,4,0,0,1)
  ,("Groenlandia",4,0,0,4,0,0,0)
  ,("Guinea",4,2,0,4,0,0,0)
  ,("Bahamas",4,0,0,4,0,0,0)
  ,("Congo",4,1,0,4,0,0,0)
  ,("Suazilandia",4,0,0,4,0,0,0)
  ,("Namibia",4,1,0,4,0,0,0)
  ,("Curazao",4,1,0,3,0,0,1)
  ,("San Bartolome",3,0,0,3,0,0,0)
  ,("Zambia",3,0,0,3,0,0,0)
  ,("Angola",3,1,0,3,0,0,0)
  ,("Republica de Yibuti",3,2,0,3,0,0,0)
  ,("Republica de Africa Central",3,0,0,3,0,0,0)
  ,("Fiyi",3,1,0,3,0,0,0)
  ,("Liberia",3,0,0,3,0,0,0)
  ,("Antigua y Barbuda",3,2,0,3,0,0,0)
  ,("Santa Lucia",3,1,0,3,0,0,0)
  ,("Cabo Verde",3,0,0,3,0,0,0)
  ,("El Salvador",3,0,0,3,0,0,0)
  ,("Zimbabue",3,0,0,2,0,1,1)
  ,("Islas Caimanes",3,0,0,2,0,0,1)
  ,("San Martin (Parte Holandesa)",2,1,0,2,0,0,0)
  ,("Niger",2,0,0,2,0,0,0)
  ,("Benin",2,0,0,2,0,0,0)
  ,("Birmania",2,2,0,2,0,0,0)
  ,("Mauritania",2,0,0,2,0,0,0)
  ,("Nicaragua",2,0,0,2,0,0,0)
  ,("Butan",2,0,0,2,0,0,0)
  ,("Sudan",2,0,0,1,0,0,1)
  ,("Gambia",2,1,0,1,0,1,1)
  ,("Nepal",2,1,1,1,0,0,0)
  ,("Belice",1,1,0,1,0,0,0)
  ,("Republica del Chad",1,0,0,1,0,0,0)
  ,("Papua Nueva Guinea",1,0,0,1,0,0,0)
  ,("Montserrat",1,0,0,1,0,0,0)
  ,("Timor Oriental",1,0,0,1,0,0,0)
  ,("Eritrea",1,0,0,1,0,0,0)
  ,("Uganda",1,0,0,1,0,0,0)
  ,("Mozambique",1,0,0,1,0,0,0)
  ,("Siria",1,0,0,1,0,0,0)
  ,("San Vicente y las Granadinas",1,0,0,1,0,0,0)
  ,("Granada",1,0,0,1,0,0,0)
  ,("Dominica",1,0,0,1,0,0,0)
  ,("Islas Turcas y Caicos",1,1,0,1,0,0,0)
  ,("Somalia",1,0,0,1,0,0,0)
  ,("Santa Sede",1,0,0,1,0,0,0)
)

$r = 4
foreach ($row in $data) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = [double]$row[1]
  $ws.Cells.Item($r, 3).Value = [double]$row[2]
  $ws.Cells.Item($r, 4).Value = [double]$row[3]
  $ws.Cells.Item($r, 5).Value = [double]$row[4]
  $ws.Cells.Item($r, 6).Value = [double]$row[5]
  $ws.Cells.Item($r, 7).Value = [double]$row[6]
  $ws.Cells.Item($r, 8).Value = [double]$row[7]
  $r = $r + 1
}
